$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$para2 = $tr.Paragraphs(2,1)
$para2.Text = "There once was a woman who was ill"
$para2b = $tr.Paragraphs(2,1)
$para2b.InsertAfter("`r")

$para3 = $tr.Paragraphs(3,1)
$para3.Text = "She kept getting worse "
$para3b = $tr.Paragraphs(3,1)
$para3b.InsertAfter("lil")
$para3c = $tr.Paragraphs(3,1)
$para3c.InsertAfter(" by ")
$para3d = $tr.Paragraphs(3,1)
$para3d.InsertAfter("lil")
$para3e = $tr.Paragraphs(3,1)
$para3e.InsertAfter("`r")

$para4 = $tr.Paragraphs(4,1)
$para4.Text = "Someone figured out what "
$para4b = $tr.Paragraphs(4,1)
$para4b.InsertAfter("was wrong,")

$para5 = $tr.Paragraphs(5,1)
$para5.Delete()
$para5b = $tr.Paragraphs(5,1)
$para5b.Delete()
